# The template's first line currently reads: {{ 'asdf' }}
# spread across three runs: "{{ ", "'asdf'", and " }}".
# Replace just the placeholder expression "'asdf'" with "today()" so the
# surrounding "{{ " / " }}" runs (and their empty <w:rPr/>) are left
# completely untouched, and the middle run becomes "today()" merging the
# whole paragraph into a single run reading "{{ today() }}".
$d = $word.ActiveDocument

$d.Content.Find.Execute("'asdf'", $false, $false, $false, $false, $false,
                         $true, 1, $false, "today()", 2)
